$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsOverview.Range("G2").Value = "2017-02-15 05:56:57"

$wsZhCn.Range("H2").Value = "2017-02-15 05:56:39"
$wsZhCn.Range("L2").Value = "2017-02-15 05:57:40"

$wsDeDe.Range("H2").Value = "2017-02-15 05:56:57"
$wsDeDe.Range("L2").Value = "2017-02-15 05:58:04"
